$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# New values keyed by actual table row number (1-based) -> array of 5 column values
$newValues = @{
    1  = @("35÷9=3, 8", "61÷5=12, 1", "15÷2=7, 1", "72÷5=14, 2", "15÷3=5, 0");
    5  = @("35÷6=5, 5", "26÷4=6, 2", "74÷7=10, 4", "27÷9=3, 0", "56÷6=9, 2");
    9  = @("18÷2=9, 0", "96÷5=19, 1", "43÷3=14, 1", "19÷8=2, 3", "87÷9=9, 6");
    13 = @("44÷3=14, 2", "59÷2=29, 1", "15÷6=2, 3", "90÷2=45, 0", "11÷6=1, 5");
    17 = @("60÷3=20, 0", "43÷4=10, 3", "72÷2=36, 0", "73÷4=18, 1", "26÷8=3, 2");
}

foreach ($r in $newValues.Keys) {
    $rowValues = $newValues[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cellRange = $cell.Range
        # Trim the trailing cell-mark / paragraph-mark character from the range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $rowValues[$c - 1]
    }
}
